$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: copy formatting for each new row from its template row ---
# (done before row 13 itself is restyled, since row 13's original B:E
#  format is the template for the '_4555' continuation-row pattern)
$ws.Range("A2:E2").Copy()
$ws.Range("A14:E14").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A3:E3").Copy()
$ws.Range("A15:E15").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A2:E2").Copy()
$ws.Range("A16:E16").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B13:E13").Copy()
$ws.Range("B17:E17").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B13:E13").Copy()
$ws.Range("B18:E18").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B13:E13").Copy()
$ws.Range("B19:E19").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B13:E13").Copy()
$ws.Range("B20:E20").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B13:E13").Copy()
$ws.Range("B21:E21").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B13:E13").Copy()
$ws.Range("B22:E22").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B13:E13").Copy()
$ws.Range("B23:E23").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B13:E13").Copy()
$ws.Range("B24:E24").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B13:E13").Copy()
$ws.Range("B25:E25").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B13:E13").Copy()
$ws.Range("B26:E26").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B13:E13").Copy()
$ws.Range("B27:E27").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B13:E13").Copy()
$ws.Range("B28:E28").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B13:E13").Copy()
$ws.Range("B29:E29").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B13:E13").Copy()
$ws.Range("B30:E30").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B13:E13").Copy()
$ws.Range("B31:E31").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B13:E13").Copy()
$ws.Range("B32:E32").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B13:E13").Copy()
$ws.Range("B33:E33").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A3:E3").Copy()
$ws.Range("A34:E34").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A4:E4").Copy()
$ws.Range("A35:E35").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A2:E2").Copy()
$ws.Range("A36:E36").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B13:E13").Copy()
$ws.Range("B37:E37").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B13:E13").Copy()
$ws.Range("B38:E38").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A3:E3").Copy()
$ws.Range("A39:E39").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("A2:E2").Copy()
$ws.Range("A40:E40").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("B13:E13").Copy()
$ws.Range("B41:E41").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# --- Step 2: fix up row 13's own formatting (border moves down onto it) ---
$ws.Range("A3:E3").Copy()
$ws.Range("A13:E13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# --- Step 3: write the text cells in the exact order the strings were
#     authored (column-block order), so the shared-string table comes
#     out in the same sequence as the target workbook ---
$ws.Cells.Item(14,"C").Value2 = ' Whoa! What kind of Pokémon is\n[CS:K]Shaymin[CR]?!'
$ws.Cells.Item(14,"A").Value2 = 'SCRIPT/P01P04A/us3105.ssb'
$ws.Cells.Item(14,"D").Value2 = ' Ого! Что за Покемоны эти\n[CS:K]Шеймины[CR]?!'
$ws.Cells.Item(14,"E").Value2 = ' Ïãï! Œóï èà Ðïëåíïîú üóé\n[CS:K]Šåêíéîú[CR]?!'
$ws.Cells.Item(15,"C").Value2 = ' So curious!'
$ws.Cells.Item(15,"D").Value2 = ' Мне так интересно! '
$ws.Cells.Item(15,"E").Value2 = ' Íîå óàë éîóåñåòîï! '
$ws.Cells.Item(16,"C").Value2 = ' Oh no!'
$ws.Cells.Item(16,"A").Value2 = 'SCRIPT/D73P21A/s32a0401.ssb'
$ws.Cells.Item(17,"C").Value2 = ' I completely forgot to bring the\n[s_item:0]!'
$ws.Cells.Item(18,"C").Value2 = ' I''m so forgetful!'
$ws.Cells.Item(19,"C").Value2 = '[CN]Give the [s_item:0]?'
$ws.Cells.Item(20,"C").Value2 = 'Yes'
$ws.Cells.Item(21,"C").Value2 = '[CN][player] gave the\n[CN][s_item:0].'
$ws.Cells.Item(22,"C").Value2 = ' Woo-hoo!\nFor me?'
$ws.Cells.Item(23,"C").Value2 = ' Really?!\nYou''re really a great friend!!'
$ws.Cells.Item(24,"C").Value2 = ' I''ll find a way to pay you back!'
$ws.Cells.Item(25,"C").Value2 = ' Oh, that''s right...! Here!'
$ws.Cells.Item(26,"C").Value2 = ' This is the best item I have\nwith me right now.[K] The [s_item:1]!'
$ws.Cells.Item(27,"C").Value2 = ' Please take this as a token of\nmy gratitude!'
$ws.Cells.Item(28,"C").Value2 = '[CN][player] received the\n[CN][s_item:1] as thanks!'
$ws.Cells.Item(29,"C").Value2 = ' I won''t forget this favor!\nThank you!'
$ws.Cells.Item(30,"C").Value2 = ' Hee-hee-hee. You''re really nice,\naren''t you, [player]?'
$ws.Cells.Item(31,"C").Value2 = ' Mountain climbing is give\nand take.'
$ws.Cells.Item(32,"C").Value2 = ' We should help others when\nthey''re in trouble.'
$ws.Cells.Item(33,"C").Value2 = '[CN]Don''t have the [s_item:0]...?'
$ws.Cells.Item(34,"C").Value2 = 'No'
$ws.Cells.Item(16,"D").Value2 = ' О, нет! '
$ws.Cells.Item(17,"D").Value2 = ' Я совершенно забыл взять с\nсобой предмет [s_item:0]!'
$ws.Cells.Item(18,"D").Value2 = ' Какой я забывчивый!'
$ws.Cells.Item(19,"D").Value2 = '[CN]Дать предмет [s_item:0]?'
$ws.Cells.Item(20,"D").Value2 = 'Да'
$ws.Cells.Item(21,"D").Value2 = '[CN][player] отдаёт предмет\n[CN][s_item:0].'
$ws.Cells.Item(22,"D").Value2 = ' Ух-ху! Это мне?'
$ws.Cells.Item(23,"D").Value2 = ' Правда?! Спасибо тебе огромное,\nдруг!!'
$ws.Cells.Item(24,"D").Value2 = ' Надо как-нибудь тебя\nотблагодарить!'
$ws.Cells.Item(25,"D").Value2 = ' О, точно!.. Вот!'
$ws.Cells.Item(26,"D").Value2 = ' Это лучшее, что я могу тебе\nдать.[K] [s_item:1]!'
$ws.Cells.Item(27,"D").Value2 = ' Прошу, прими это в знак\nпризнательности!'
$ws.Cells.Item(28,"D").Value2 = '[CN][player] получает предмет\n[CN][s_item:1]\n[CN]в знак признательности!'
$ws.Cells.Item(29,"D").Value2 = ' Я этого никогда не забуду!\nСпасибо тебе!'
$ws.Cells.Item(30,"D").Value2 = ' Хи-хи-хи. Не можешь пройти мимо\nбеды, да, [player]?'
$ws.Cells.Item(31,"D").Value2 = ' Восхождение на гору даётся\nнепросто.'
$ws.Cells.Item(32,"D").Value2 = ' Нам нужно помогать тем, кто\nнуждается в помощи.'
$ws.Cells.Item(33,"D").Value2 = '[CN]Нет предмета [s_item:0]?..'
$ws.Cells.Item(34,"D").Value2 = 'Нет'
$ws.Cells.Item(16,"E").Value2 = ' Ï, îåó! '
$ws.Cells.Item(17,"E").Value2 = ' Ÿ òïâåñšåîîï èàáúì âèÿóû ò\nòïáïê ðñåäíåó [s_item:0]!'
$ws.Cells.Item(18,"E").Value2 = ' Ëàëïê ÿ èàáúâœéâúê!'
$ws.Cells.Item(19,"E").Value2 = '[CN]Äàóû ðñåäíåó [s_item:0]?'
$ws.Cells.Item(20,"E").Value2 = 'Äà'
$ws.Cells.Item(21,"E").Value2 = '[CN][player] ïóäàæó ðñåäíåó\n[CN][s_item:0].'
$ws.Cells.Item(22,"E").Value2 = ' Ôö-öô! Üóï íîå?'
$ws.Cells.Item(23,"E").Value2 = ' Ðñàâäà?! Òðàòéáï óåáå ïãñïíîïå,\näñôã!!'
$ws.Cells.Item(24,"E").Value2 = ' Îàäï ëàë-îéáôäû óåáÿ\nïóáìàãïäàñéóû!'
$ws.Cells.Item(25,"E").Value2 = ' Ï, óïœîï!.. Âïó!'
$ws.Cells.Item(26,"E").Value2 = ' Üóï ìôœšåå, œóï ÿ íïãô óåáå\näàóû.[K] [s_item:1]!'
$ws.Cells.Item(27,"E").Value2 = ' Ðñïšô, ðñéíé üóï â èîàë\nðñéèîàóåìûîïòóé!'
$ws.Cells.Item(28,"E").Value2 = '[CN][player] ðïìôœàåó ðñåäíåó\n[CN][s_item:1]\n[CN]â èîàë ðñéèîàóåìûîïòóé!'
$ws.Cells.Item(29,"E").Value2 = ' Ÿ üóïãï îéëïãäà îå èàáôäô!\nÒðàòéáï óåáå!'
$ws.Cells.Item(30,"E").Value2 = ' Öé-öé-öé. Îå íïçåšû ðñïêóé íéíï\náåäú, äà, [player]?'
$ws.Cells.Item(31,"E").Value2 = ' Âïòöïçäåîéå îà ãïñô äàæóòÿ\nîåðñïòóï.'
$ws.Cells.Item(32,"E").Value2 = ' Îàí îôçîï ðïíïãàóû óåí, ëóï\nîôçäàåóòÿ â ðïíïþé.'
$ws.Cells.Item(33,"E").Value2 = '[CN]Îåó ðñåäíåóà [s_item:0]?..'
$ws.Cells.Item(34,"E").Value2 = 'Îåó'
$ws.Cells.Item(35,"C").Value2 = ' I will never forget yoooou!\nThank yoooou!!'
$ws.Cells.Item(35,"A").Value2 = 'SCRIPT/D73P21A/us3106.ssb'
$ws.Cells.Item(35,"D").Value2 = ' Я тебяяяя никогдаааа не\nзабууууду! Спасиииибо!!'
$ws.Cells.Item(35,"E").Value2 = ' Ÿ óåáÿÿÿÿ îéëïãäàààà îå\nèàáôôôôäô! Òðàòééééáï!!'
$ws.Cells.Item(36,"C").Value2 = ' Whoa! Still the 3rd Station\nCleeeearing!!'
$ws.Cells.Item(37,"C").Value2 = ' I heard there''s a gondola to\nthe 6th Station Clearing…'
$ws.Cells.Item(38,"C").Value2 = ' But! I want to climb\nthere myself!!'
$ws.Cells.Item(39,"C").Value2 = ' Whoa!\nWhat willpower!!'
$ws.Cells.Item(36,"A").Value2 = 'SCRIPT/D73P23A/us3108.ssb'
$ws.Cells.Item(36,"D").Value2 = ' Вау! Я на Поляяяяне 3-го\nПерехооода!!'
$ws.Cells.Item(37,"D").Value2 = ' Говорят, на Поляне 6-го\nПерехода есть Гондола...'
$ws.Cells.Item(38,"D").Value2 = ' Но! Я сам хочу туда дойти!!'
$ws.Cells.Item(39,"D").Value2 = ' Вау! Вот это у меня сила воли!!'
$ws.Cells.Item(36,"E").Value2 = ' Âàô! Ÿ îà Ðïìÿÿÿÿîå 3-ãï\nÐåñåöïïïäà!!'
$ws.Cells.Item(37,"E").Value2 = ' Ãïâïñÿó, îà Ðïìÿîå 6-ãï\nÐåñåöïäà åòóû Ãïîäïìà...'
$ws.Cells.Item(38,"E").Value2 = ' Îï! Ÿ òàí öïœô óôäà äïêóé!!'
$ws.Cells.Item(39,"E").Value2 = ' Âàô! Âïó üóï ô íåîÿ òéìà âïìé!!'
$ws.Cells.Item(40,"C").Value2 = ' Whoa!\nWe caught up to [CS:N]Mr. Mime[CR]''s team!'
$ws.Cells.Item(41,"C").Value2 = ' This makes me haaaappy!!'
$ws.Cells.Item(40,"A").Value2 = 'SCRIPT/D73P27A/us3107.ssb'
$ws.Cells.Item(40,"D").Value2 = ' Вау! Мы догнали команду\n[CS:N]Мр-а Майма[CR]!'
$ws.Cells.Item(41,"D").Value2 = ' Я так счааааастлив!!'
$ws.Cells.Item(40,"E").Value2 = ' Âàô! Íú äïãîàìé ëïíàîäô\n[CS:N]Íñ-à Íàêíà[CR]!'
$ws.Cells.Item(41,"E").Value2 = ' Ÿ óàë òœàààààòóìéâ!!'

# --- Step 4: write the numeric 'line number' column (B) ---
$ws.Cells.Item(13,2).Value2 = 229
$ws.Cells.Item(14,2).Value2 = 192
$ws.Cells.Item(15,2).Value2 = 202
$ws.Cells.Item(16,2).Value2 = 17
$ws.Cells.Item(17,2).Value2 = 20
$ws.Cells.Item(18,2).Value2 = 27
$ws.Cells.Item(19,2).Value2 = 32
$ws.Cells.Item(20,2).Value2 = 36
$ws.Cells.Item(21,2).Value2 = 47
$ws.Cells.Item(22,2).Value2 = 66
$ws.Cells.Item(23,2).Value2 = 73
$ws.Cells.Item(24,2).Value2 = 76
$ws.Cells.Item(25,2).Value2 = 92
$ws.Cells.Item(26,2).Value2 = 95
$ws.Cells.Item(27,2).Value2 = 98
$ws.Cells.Item(28,2).Value2 = 103
$ws.Cells.Item(29,2).Value2 = 113
$ws.Cells.Item(30,2).Value2 = 133
$ws.Cells.Item(31,2).Value2 = 145
$ws.Cells.Item(32,2).Value2 = 148
$ws.Cells.Item(33,2).Value2 = 164
$ws.Cells.Item(34,2).Value2 = 171
$ws.Cells.Item(35,2).Value2 = 163
$ws.Cells.Item(36,2).Value2 = 129
$ws.Cells.Item(37,2).Value2 = 136
$ws.Cells.Item(38,2).Value2 = 139
$ws.Cells.Item(39,2).Value2 = 146
$ws.Cells.Item(40,2).Value2 = 109
$ws.Cells.Item(41,2).Value2 = 112

# --- Step 5: explicit row heights (rows whose wrapped content needs
#     more than the default single-line height) ---
$ws.Rows.Item(14).RowHeight = 43.2
$ws.Rows.Item(16).RowHeight = 43.2
$ws.Rows.Item(17).RowHeight = 21.6
$ws.Rows.Item(21).RowHeight = 21.6
$ws.Rows.Item(23).RowHeight = 22.8
$ws.Rows.Item(24).RowHeight = 21.6
$ws.Rows.Item(26).RowHeight = 21.6
$ws.Rows.Item(27).RowHeight = 21.6
$ws.Rows.Item(28).RowHeight = 31.8
$ws.Rows.Item(29).RowHeight = 27
$ws.Rows.Item(30).RowHeight = 21.6
$ws.Rows.Item(31).RowHeight = 21.6
$ws.Rows.Item(32).RowHeight = 21.6
$ws.Rows.Item(35).RowHeight = 43.2
$ws.Rows.Item(36).RowHeight = 43.2
$ws.Rows.Item(37).RowHeight = 21.6
$ws.Rows.Item(40).RowHeight = 43.2

# --- Step 6: selection / scroll position to match the end of the edit ---
$ws.Range("D40").Select()
